# Auto update Excel log
# Appends newly-collected sensor readings to the PIR, Humidity and
# Temperature sheets (mirrors the source system's log-append job).

$wb = $excel.ActiveWorkbook

function Add-LogRows($SheetName, $StartRow, $Rows, $TextColumns) {
    $ws = $wb.Worksheets.Item($SheetName)

    for ($i = 0; $i -lt $Rows.Length; $i++) {
        $r = $StartRow + $i
        $rowData = $Rows[$i]

        for ($c = 0; $c -lt $rowData.Length; $c++) {
            $col = $c + 1
            $cell = $ws.Cells.Item($r, $col)

            if ($TextColumns -contains $col) {
                # Force plain text so values like dates / percentages
                # are stored verbatim instead of being auto-converted
                # into numeric date / percentage values.
                $cell.NumberFormat = "@"
            }

            $cell.Value = $rowData[$c]
        }
    }
}

# --- PIR: rows 65-77 (Bathroom / No Motion / Inactive) ---------------------
$pirRows = @(
    @("2026-01-30","15:42:44","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","15:42:48","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","15:42:53","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","15:42:58","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","15:43:03","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","15:43:08","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","15:43:13","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","15:43:18","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","15:43:23","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","15:43:28","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","15:43:33","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","15:43:38","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","15:43:43","15:00","Bathroom","No Motion","Inactive")
)
Add-LogRows "PIR" 65 $pirRows @(1)

# --- Humidity: rows 35-42 (Bathroom / humidity % / Active) -----------------
$humidityRows = @(
    @("2026-01-30","15:42:44","15:00","Bathroom","86.9%","Active"),
    @("2026-01-30","15:42:48","15:00","Bathroom","86.4%","Active"),
    @("2026-01-30","15:43:08","15:00","Bathroom","87.8%","Active"),
    @("2026-01-30","15:43:13","15:00","Bathroom","87.9%","Active"),
    @("2026-01-30","15:43:23","15:00","Bathroom","87.8%","Active"),
    @("2026-01-30","15:43:28","15:00","Bathroom","87.9%","Active"),
    @("2026-01-30","15:43:33","15:00","Bathroom","87.9%","Active"),
    @("2026-01-30","15:43:43","15:00","Bathroom","87.9%","Active")
)
Add-LogRows "Humidity" 35 $humidityRows @(1, 5)

# --- Temperature: rows 2-5 (Living Room / mmWave presence / status) --------
$temperatureRows = @(
    @("2026-01-30","15:43:03","15:00","Living Room","NO_MOTION_DETECTED","Inactive"),
    @("2026-01-30","15:43:14","15:00","Living Room","PRESENCE_DETECTED","Active"),
    @("2026-01-30","15:43:24","15:00","Living Room","PRESENCE_DETECTED","Active"),
    @("2026-01-30","15:43:35","15:00","Living Room","PRESENCE_DETECTED","Active")
)
Add-LogRows "Temperature" 2 $temperatureRows @(1)
